# Append the new daily mods-count row (2026/01/16) to the tracking sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 66
$newRow  = $lastRow + 1

# --- Column A: the date, written as literal text -----------------------
# A plain $ws.Cells.Item($newRow,1).Value = "2026/01/16" would be
# auto-recognized as a date (like real Excel does for unformatted cells)
# and stored as a serial number, which is not what the sheet uses for its
# other rows (they store the date as literal text). To keep it a literal
# string we stage the text as a formula result in a scratch cell (so it is
# never subject to Excel's "looks like a date" input heuristic) and copy
# only the *value* over with PasteSpecial.
$scratch = $ws.Range("Z1")
$scratch.Formula = "=""2026/01/16"""
$scratch.Copy()
$ws.Range("A" + $newRow).PasteSpecial(-4163)   # xlPasteValues
$scratch.ClearContents()

# --- Columns B and C: plain text / number, no special parsing needed ---
$ws.Cells.Item($newRow, 2).Value = "逃离鸭科夫"
$ws.Cells.Item($newRow, 3).Value = 1146

# --- Formatting: mirror the style used by the rest of the data rows ----
$ws.Range("A" + $lastRow + ":C" + $lastRow).Copy()
$ws.Range("A" + $newRow + ":C" + $newRow).PasteSpecial(-4122)   # xlPasteFormats
